$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blue fill used for "filled in" hour cells (theme color 0070C0 -> BGR 12611584)
$blue = 12611584
# Red fill used for the zero/"invalid" hour cell in column D (FF0000 -> BGR 255)
$red = 255

# Row 101 (name "Ellen"-ish / 4th person row in the weekly block) - fill in hours
$ws.Range("B101").Value = 4
$ws.Range("C101").Value = 4
$ws.Range("C101").Interior.Color = $blue
$ws.Range("D101").Value = 0
$ws.Range("D101").Interior.Color = $red
$ws.Range("E101").Value = 4
$ws.Range("E101").Interior.Color = $blue
$ws.Range("F101").Value = 4
$ws.Range("F101").Interior.Color = $blue
$ws.Range("G101").Value = 4
$ws.Range("G101").Interior.Color = $blue

# Row 102 (5th person row in the weekly block) - fill in hours
$ws.Range("B102").Value = 4
$ws.Range("C102").Value = 4
$ws.Range("C102").Interior.Color = $blue
$ws.Range("D102").Value = 0
$ws.Range("D102").Interior.Color = $red
$ws.Range("E102").Value = 4
$ws.Range("E102").Interior.Color = $blue
$ws.Range("F102").Value = 4
$ws.Range("F102").Interior.Color = $blue
$ws.Range("G102").Value = 4
$ws.Range("G102").Interior.Color = $blue

# Scroll back up and select I20 instead of I97 (matches the saved view state)
$ws.Range("I20").Select() | Out-Null
